$wb = $excel.ActiveWorkbook

# Rename customer header columns: first_name -> f_name, last_name -> l_name
$wsCustomers = $wb.Worksheets.Item("customers")
$wsCustomers.Range("B1").Value = "f_name"
$wsCustomers.Range("C1").Value = "l_name"

# Move the cursor on the instructors sheet (non-active tab) first so the
# later selection on "customers" below is the one that ends up marking the
# active sheet / tab.
$wsInstructors = $wb.Worksheets.Item("instructors")
$wsInstructors.Range("C3").Select()

# customers sheet stays the active tab with the selection moved to D38
$wsCustomers.Range("D38").Select()
